$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (B2:E2, G2) - F2 stays 0
$ws.Range("B2").Value = 0.1190320826869504
$ws.Range("C2").Value = 0.306821227259698
$ws.Range("D2").Value = 0.7527432677738641
$ws.Range("E2").Value = 10.19245300693656
$ws.Range("G2").Value = 11.37104958465707

# Row 3 updates (B3:E3, G3) - F3 stays 1
$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 10.34677158129881
$ws.Range("D3").Value = 22.3905356188092
$ws.Range("E3").Value = 1133.036916526867
$ws.Range("G3").Value = 1169.06105627184
